$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Sheet shuffle: turn the single "2022-Q1" sheet into two sheets,
#    "2022-Q3" (new data, keeps the original sheetId/rId so it lands right
#    after "总计") and "2022-Q1" (a duplicate holding the untouched old data).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2022-Q1").Name = "2022-Q3"

$srcQ3 = $wb.Worksheets.Item("2022-Q3")
$srcQ3.Copy($null, $srcQ3)
$wb.Worksheets.Item("2022-Q3 (2)").Name = "2022-Q1"

# ---------------------------------------------------------------------------
# 2. "总计" (summary) sheet: insert the 2022-Q3 row above the 2022-Q1 row.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# carry the header/"index-column" style (s=2) down onto the new row as well
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A2:A3").PasteSpecial(-4122)

$wsTotal.Range("A3").Value2 = 1
$wsTotal.Range("B3").Value2 = "2022-Q1"
$wsTotal.Range("C3").Value2 = 2
$wsTotal.Range("D3").Value2 = 0

$wsTotal.Range("A2").Value2 = 0
$wsTotal.Range("B2").Value2 = "2022-Q3"
$wsTotal.Range("C2").Value2 = 3
$wsTotal.Range("D2").Value2 = 0.08

# ---------------------------------------------------------------------------
# 3. "2022-Q3" sheet: replace the (copied) old fund table with the new one.
# ---------------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Item("2022-Q3")

# Column codes / ratios are text in the source data (leading zeros etc.) -
# force text storage with a quote-prefix, then drop the resulting format so
# no stray number-format style is left behind.
$wsQ3.Range("B2:G4").Value2 = "'"
$wsQ3.Range("B2").Value2 = "'016029"
$wsQ3.Range("C2").Value2 = "湘财成长优选一年持有期混合A"
$wsQ3.Range("D2").Value2 = "'1.77"
$wsQ3.Range("E2").Value2 = "'75.14"
$wsQ3.Range("F2").Value2 = "'4.20"
$wsQ3.Range("G2").Value2 = "'0.0743"
$wsQ3.Range("H2").Value2 = 2

$wsQ3.Range("B3").Value2 = "'016030"
$wsQ3.Range("C3").Value2 = "湘财成长优选一年持有期混合C"
$wsQ3.Range("D3").Value2 = "'0.16"
$wsQ3.Range("E3").Value2 = "'75.14"
$wsQ3.Range("F3").Value2 = "'4.20"
$wsQ3.Range("G3").Value2 = "'0.0067"
$wsQ3.Range("H3").Value2 = 2

$wsQ3.Range("B4").Value2 = "'002000"
$wsQ3.Range("C4").Value2 = "工银新生利混合"
$wsQ3.Range("D4").Value2 = "'1.08"
$wsQ3.Range("E4").Value2 = "'28.83"
$wsQ3.Range("F4").Value2 = "'0.08"
$wsQ3.Range("G4").Value2 = "'0.0009"
$wsQ3.Range("H4").Value2 = 10

$wsQ3.Range("A2").Value2 = 0
$wsQ3.Range("A3").Value2 = 1
$wsQ3.Range("A4").Value2 = 2

$wsQ3.Range("B2:G4").ClearFormats()

# Headers (row 1) - same labels, just re-typed to be explicit.
$wsQ3.Range("B1").Value2 = "基金代码"
$wsQ3.Range("C1").Value2 = "基金名称"
$wsQ3.Range("D1").Value2 = "基金规模"
$wsQ3.Range("E1").Value2 = "股票总仓位"
$wsQ3.Range("F1").Value2 = "仓位占比"
$wsQ3.Range("G1").Value2 = "持有市值(亿元)"
$wsQ3.Range("H1").Value2 = "仓位排名"

# Re-apply the bold/bordered "总计"-style header formatting (style index 2)
# to the header row and the index column, matching the rest of the workbook.
$wsTotal.Range("B1").Copy()
$wsQ3.Range("B1:H1").PasteSpecial(-4122)
$wsTotal.Range("A2").Copy()
$wsQ3.Range("A2:A4").PasteSpecial(-4122)
